$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Fgf7"
$ws.Cells.Item(2, 3).Value = "Fgfr2"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.245245
$ws.Cells.Item(2, 8).Value = 0.735735
$ws.Cells.Item(2, 9).Value = 0.0130094690177091
$ws.Cells.Item(2, 10).Value = 0.0130094690177091
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.09434
$ws.Cells.Item(2, 14).Value = 0.28302
$ws.Cells.Item(2, 15).Value = 0.05191071108246543
$ws.Cells.Item(2, 16).Value = 0.05191071108246543
$ws.Cells.Item(2, 17).Value = 0.0231364133
$ws.Cells.Item(2, 18).Value = 0.2082277197
$ws.Cells.Item(2, 19).Value = 0.0006753307875145823
$ws.Cells.Item(2, 20).Value = 0.0006753307875145823

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Fgf7"
$ws.Cells.Item(3, 3).Value = "Fgfr2"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.245245
$ws.Cells.Item(3, 8).Value = 0.735735
$ws.Cells.Item(3, 9).Value = 0.0130094690177091
$ws.Cells.Item(3, 10).Value = 0.0130094690177091
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.9431116666666667
$ws.Cells.Item(3, 14).Value = 2.829335
$ws.Cells.Item(3, 15).Value = 0.5189484550226392
$ws.Cells.Item(3, 16).Value = 0.5189484550226391
$ws.Cells.Item(3, 17).Value = 0.2312934206916667
$ws.Cells.Item(3, 18).Value = 2.081640786225
$ws.Cells.Item(3, 19).Value = 0.006751243847405029
$ws.Cells.Item(3, 20).Value = 0.006751243847405027

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Fgf7"
$ws.Cells.Item(4, 3).Value = "Fgfr2"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.245245
$ws.Cells.Item(4, 8).Value = 0.735735
$ws.Cells.Item(4, 9).Value = 0.0130094690177091
$ws.Cells.Item(4, 10).Value = 0.0130094690177091
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.7798996666666667
$ws.Cells.Item(4, 14).Value = 2.339699
$ws.Cells.Item(4, 15).Value = 0.4291408338948954
$ws.Cells.Item(4, 16).Value = 0.4291408338948954
$ws.Cells.Item(4, 17).Value = 0.1912664937516667
$ws.Cells.Item(4, 18).Value = 1.721398443765
$ws.Cells.Item(4, 19).Value = 0.005582894382789488
$ws.Cells.Item(4, 20).Value = 0.005582894382789488

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Fgf7"
$ws.Cells.Item(5, 3).Value = "Fgfr2"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 17.39906333333333
$ws.Cells.Item(5, 8).Value = 52.19719000000001
$ws.Cells.Item(5, 9).Value = 0.9229650976458578
$ws.Cells.Item(5, 10).Value = 0.9229650976458579
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.09434
$ws.Cells.Item(5, 14).Value = 0.28302
$ws.Cells.Item(5, 15).Value = 0.05191071108246543
$ws.Cells.Item(5, 16).Value = 0.05191071108246543
$ws.Cells.Item(5, 17).Value = 1.641427634866667
$ws.Cells.Item(5, 18).Value = 14.7728487138
$ws.Cells.Item(5, 19).Value = 0.04791177452309361
$ws.Cells.Item(5, 20).Value = 0.04791177452309362

# Row 6
$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Fgf7"
$ws.Cells.Item(6, 3).Value = "Fgfr2"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 17.39906333333333
$ws.Cells.Item(6, 8).Value = 52.19719000000001
$ws.Cells.Item(6, 9).Value = 0.9229650976458578
$ws.Cells.Item(6, 10).Value = 0.9229650976458579
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.9431116666666667
$ws.Cells.Item(6, 14).Value = 2.829335
$ws.Cells.Item(6, 15).Value = 0.5189484550226392
$ws.Cells.Item(6, 16).Value = 0.5189484550226391
$ws.Cells.Item(6, 17).Value = 16.40925961873889
$ws.Cells.Item(6, 18).Value = 147.68333656865
$ws.Cells.Item(6, 19).Value = 0.4789713114631372
$ws.Cells.Item(6, 20).Value = 0.4789713114631372

# Row 7
$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Fgf7"
$ws.Cells.Item(7, 3).Value = "Fgfr2"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 17.39906333333333
$ws.Cells.Item(7, 8).Value = 52.19719000000001
$ws.Cells.Item(7, 9).Value = 0.9229650976458578
$ws.Cells.Item(7, 10).Value = 0.9229650976458579
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.7798996666666667
$ws.Cells.Item(7, 14).Value = 2.339699
$ws.Cells.Item(7, 15).Value = 0.4291408338948954
$ws.Cells.Item(7, 16).Value = 0.4291408338948954
$ws.Cells.Item(7, 17).Value = 13.56952369397889
$ws.Cells.Item(7, 18).Value = 122.12571324581
$ws.Cells.Item(7, 19).Value = 0.3960820116596269
$ws.Cells.Item(7, 20).Value = 0.396082011659627

# Row 8
$ws.Cells.Item(8, 1).Value = "ECs"
$ws.Cells.Item(8, 2).Value = "Fgf7"
$ws.Cells.Item(8, 3).Value = "Fgfr2"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.206960666666667
$ws.Cells.Item(8, 8).Value = 3.620882
$ws.Cells.Item(8, 9).Value = 0.06402543333643303
$ws.Cells.Item(8, 10).Value = 0.06402543333643303
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.09434
$ws.Cells.Item(8, 14).Value = 0.28302
$ws.Cells.Item(8, 15).Value = 0.05191071108246543
$ws.Cells.Item(8, 16).Value = 0.05191071108246543
$ws.Cells.Item(8, 17).Value = 0.1138646692933333
$ws.Cells.Item(8, 18).Value = 1.02478202364
$ws.Cells.Item(8, 19).Value = 0.003323605771857226
$ws.Cells.Item(8, 20).Value = 0.003323605771857226

# Row 9
$ws.Cells.Item(9, 1).Value = "ECs"
$ws.Cells.Item(9, 2).Value = "Fgf7"
$ws.Cells.Item(9, 3).Value = "Fgfr2"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.206960666666667
$ws.Cells.Item(9, 8).Value = 3.620882
$ws.Cells.Item(9, 9).Value = 0.06402543333643303
$ws.Cells.Item(9, 10).Value = 0.06402543333643303
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.9431116666666667
$ws.Cells.Item(9, 14).Value = 2.829335
$ws.Cells.Item(9, 15).Value = 0.5189484550226392
$ws.Cells.Item(9, 16).Value = 0.5189484550226391
$ws.Cells.Item(9, 17).Value = 1.138298685941111
$ws.Cells.Item(9, 18).Value = 10.24468817347
$ws.Cells.Item(9, 19).Value = 0.0332258997120969
$ws.Cells.Item(9, 20).Value = 0.03322589971209689

# Row 10
$ws.Cells.Item(10, 1).Value = "ECs"
$ws.Cells.Item(10, 2).Value = "Fgf7"
$ws.Cells.Item(10, 3).Value = "Fgfr2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.206960666666667
$ws.Cells.Item(10, 8).Value = 3.620882
$ws.Cells.Item(10, 9).Value = 0.06402543333643303
$ws.Cells.Item(10, 10).Value = 0.06402543333643303
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.7798996666666667
$ws.Cells.Item(10, 14).Value = 2.339699
$ws.Cells.Item(10, 15).Value = 0.4291408338948954
$ws.Cells.Item(10, 16).Value = 0.4291408338948954
$ws.Cells.Item(10, 17).Value = 0.941308221613111
$ws.Cells.Item(10, 18).Value = 8.471773994518
$ws.Cells.Item(10, 19).Value = 0.0274759278524789
$ws.Cells.Item(10, 20).Value = 0.0274759278524789
